# Auto-generated edit script applying the cryptos.xlsx diff
# (update of Price (D) and Volume/1h (E) columns, plus a couple of row
# content swaps for rows 30/31 and 44/45).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '92.603.22'
$ws.Range('E2').Value = '  -1.86%  '
$ws.Range('D3').Value = '3.406.33'
$ws.Range('E3').Value = '  -0.44%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').Value = "'230.47"
$ws.Range('E5').Value = '  -2.91%  '
$ws.Range('D6').Value = "'617.01"
$ws.Range('E6').Value = '  -4.01%  '
$ws.Range('D7').Value = "'1.37"
$ws.Range('E7').Value = '  -5.08%  '
$ws.Range('D8').Value = "'0.389"
$ws.Range('E8').Value = '  -4.02%  '
$ws.Range('E9').Value = '  +0.09%  '
$ws.Range('D10').Value = "'0.960"
$ws.Range('E10').Value = '  -0.84%  '
$ws.Range('D11').Value = '3.403.70'
$ws.Range('D12').Value = "'42.87"
$ws.Range('E12').Value = '  +3.34%  '
$ws.Range('E13').Value = '  -1.04%  '
$ws.Range('D14').Value = "'6.20"
$ws.Range('E14').Value = '  -0.10%  '
$ws.Range('D15').Value = '4.054.57'
$ws.Range('E15').Value = '  -0.25%  '
$ws.Range('D16').Value = '92.612.60'
$ws.Range('E16').Value = '  -1.62%  '
$ws.Range('D17').Value = "'0.0000245"
$ws.Range('E17').Value = '  -2.68%  '
$ws.Range('D18').Value = "'8.05"
$ws.Range('E18').Value = '  -3.00%  '
$ws.Range('D19').Value = '3.410.22'
$ws.Range('E19').Value = '  -0.22%  '
$ws.Range('D20').Value = "'17.84"
$ws.Range('E20').Value = '  +2.21%  '
$ws.Range('D21').Value = "'11.47"
$ws.Range('E21').Value = '  -0.51%  '
$ws.Range('D22').Value = "'495.16"
$ws.Range('E22').Value = '  -0.51%  '
$ws.Range('D23').Value = "'3.30"
$ws.Range('E23').Value = '  +2.29%  '
$ws.Range('D24').Value = "'0.435"
$ws.Range('E24').Value = '  -12.72%  '
$ws.Range('D25').Value = "'6.50"
$ws.Range('E25').Value = '  +0.26%  '
$ws.Range('E26').Value = '  -5.20%  '
$ws.Range('D27').Value = "'91.39"
$ws.Range('E27').Value = '  +0.08%  '
$ws.Range('D28').Value = "'11.92"
$ws.Range('E28').Value = '  -0.17%  '
$ws.Range('D29').Value = '3.592.03'
$ws.Range('E29').Value = '  -0.38%  '
$ws.Range('B30').Value = 'Dai'
$ws.Range('C30').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D30').Value = "'1.00"
$ws.Range('E30').Value = '  +0.03%  '
$ws.Range('B31').Value = 'InternetComputer(DFINITY)'
$ws.Range('C31').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D31').Value = "'11.22"
$ws.Range('E31').Value = '  -3.72%  '
$ws.Range('D32').Value = "'2.70"
$ws.Range('E32').Value = '  -2.24%  '
$ws.Range('D33').Value = "'0.134"
$ws.Range('E33').Value = '  -3.82%  '
$ws.Range('D34').Value = "'1.00"
$ws.Range('E34').Value = '  +0.53%  '
$ws.Range('E35').Value = '  -3.89%  '
$ws.Range('D36').Value = "'29.36"
$ws.Range('E36').Value = '  -1.06%  '
$ws.Range('D37').Value = "'0.538"
$ws.Range('E37').Value = '  -2.47%  '
$ws.Range('D38').Value = "'551.51"
$ws.Range('E38').Value = '  +0.34%  '
$ws.Range('D39').Value = "'7.44"
$ws.Range('E39').Value = '  -2.61%  '
$ws.Range('E40').Value = '  -0.04%  '
$ws.Range('D41').Value = "'0.149"
$ws.Range('E41').Value = '  -0.94%  '
$ws.Range('D42').Value = "'1.38"
$ws.Range('E42').Value = '  -4.94%  '
$ws.Range('D43').Value = "'0.907"
$ws.Range('E43').Value = '  +0.69%  '
$ws.Range('B44').Value = 'MantraDAO'
$ws.Range('C44').Value = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
$ws.Range('D44').Value = "'3.72"
$ws.Range('E44').Value = '  +2.50%  '
$ws.Range('B45').Value = 'ImmutableX'
$ws.Range('C45').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D45').Value = "'1.71"
$ws.Range('E45').Value = '  -0.64%  '
$ws.Range('D46').Value = "'23.64"
$ws.Range('E46').Value = '  -1.78%  '
$ws.Range('D47').Value = "'5.45"
$ws.Range('E47').Value = '  -2.58%  '
$ws.Range('D48').Value = "'0.0404"
$ws.Range('E48').Value = '  -1.14%  '
$ws.Range('D49').Value = "'53.00"
$ws.Range('E49').Value = '  -3.39%  '
$ws.Range('D50').Value = "'2.09"
$ws.Range('E50').Value = '  -4.41%  '
$ws.Range('E51').Value = '  +14.86%  '
